$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.372.82"
$ws.Range("E2").Value = "'  -0.21%  "
$ws.Range("D3").Value = "'1.842.55"
$ws.Range("D5").Value = "'239.02"
$ws.Range("E5").Value = "'  -0.70%  "
$ws.Range("D6").Value = "'0.6310"
$ws.Range("E6").Value = "'  -0.25%  "
$ws.Range("E7").Value = "'  +0.03%  "
$ws.Range("D8").Value = "'0.07529"
$ws.Range("E8").Value = "'  -0.85%  "
$ws.Range("D9").Value = "'0.2929"
$ws.Range("E9").Value = "'  -1.43%  "
$ws.Range("D10").Value = "'24.45"
$ws.Range("E10").Value = "'  -0.68%  "
$ws.Range("D11").Value = "'0.07706"
$ws.Range("E11").Value = "'  -0.10%  "
$ws.Range("D12").Value = "'1.855.49"
$ws.Range("E12").Value = "'  -6.53%  "
$ws.Range("D13").Value = "'4.991"
$ws.Range("E13").Value = "'  +0.04%  "
$ws.Range("D14").Value = "'0.6782"
$ws.Range("E14").Value = "'  -1.05%  "
$ws.Range("D15").Value = "'0.00001035"
$ws.Range("E15").Value = "'  +3.39%  "
$ws.Range("D16").Value = "'82.81"
$ws.Range("E16").Value = "'  -0.07%  "
$ws.Range("D17").Value = "'2.122.15"
$ws.Range("E17").Value = "'  -6.28%  "
$ws.Range("D18").Value = "'6.131"
$ws.Range("E18").Value = "'  -0.98%  "
$ws.Range("D19").Value = "'29.411.42"
$ws.Range("E19").Value = "'  -0.13%  "
$ws.Range("D20").Value = "'227.84"
$ws.Range("E20").Value = "'  -1.94%  "
$ws.Range("E21").Value = "'  -0.79%  "
$ws.Range("D22").Value = "'0.9996"
$ws.Range("E22").Value = "'  -0.02%  "
$ws.Range("D23").Value = "'7.435"
$ws.Range("E23").Value = "'  -1.93%  "
$ws.Range("D24").Value = "'1.001"
$ws.Range("E24").Value = "'  +0.13%  "
$ws.Range("D25").Value = "'156.83"
$ws.Range("E25").Value = "'  +1.18%  "
$ws.Range("D26").Value = "'0.1390"
$ws.Range("E26").Value = "'  -0.05%  "
$ws.Range("E27").Value = "'  -1.13%  "
$ws.Range("E28").Value = "'  -0.49%  "
$ws.Range("D29").Value = "'1.455"
$ws.Range("E29").Value = "'  -1.18%  "
$ws.Range("D30").Value = "'1.275"
$ws.Range("E30").Value = "'  +1.22%  "
$ws.Range("D31").Value = "'0.05621"
$ws.Range("E31").Value = "'  -3.16%  "
$ws.Range("D32").Value = "'4.103"
$ws.Range("E32").Value = "'  -0.53%  "
$ws.Range("E33").Value = "'  -0.19%  "
$ws.Range("D34").Value = "'1.830"
$ws.Range("E34").Value = "'  -2.06%  "
$ws.Range("D35").Value = "'1.156"
$ws.Range("E35").Value = "'  -0.29%  "
$ws.Range("D36").Value = "'0.7064"
$ws.Range("E36").Value = "'  -2.22%  "
$ws.Range("E37").Value = "'  -0.29%  "
$ws.Range("D38").Value = "'1.240.70"
$ws.Range("E38").Value = "'  -0.74%  "
$ws.Range("E39").Value = "'  -0.14%  "
$ws.Range("D40").Value = "'2.763"
$ws.Range("E40").Value = "'  -1.09%  "
$ws.Range("D41").Value = "'6.290"
$ws.Range("E41").Value = "'  +3.24%  "
$ws.Range("D42").Value = "'0.8999"
$ws.Range("E42").Value = "'  -0.31%  "
$ws.Range("D43").Value = "'0.9992"
$ws.Range("D44").Value = "'101.90"
$ws.Range("E44").Value = "'  +0.44%  "
$ws.Range("D45").Value = "'65.49"
$ws.Range("E45").Value = "'  -2.19%  "
$ws.Range("D46").Value = "'0.00000000119"
$ws.Range("E46").Value = "'  +0.93%  "
$ws.Range("D47").Value = "'7.052"
$ws.Range("E47").Value = "'  -3.81%  "
$ws.Range("D48").Value = "'0.3997"
$ws.Range("E48").Value = "'  -0.54%  "
$ws.Range("B49").Value = "'EnergySwap"
$ws.Range("C49").Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'8.883"
$ws.Range("E49").Value = "'  -3.41%  "
$ws.Range("B50").Value = "'RenderToken"
$ws.Range("C50").Value = "'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D50").Value = "'1.663"
$ws.Range("E50").Value = "'  -1.90%  "
$ws.Range("D51").Value = "'0.1117"
$ws.Range("E51").Value = "'  -0.73%  "
